$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.868.67"
$ws.Range("E2").Value = "  +0.47%  "
$ws.Range("D3").Value = "3.494.40"
$ws.Range("E3").Value = "  +0.25%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.51"
$ws.Range("E5").Value = "  +0.53%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.60"
$ws.Range("E6").Value = "  +1.71%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.590"
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("E9").Value = "  +3.85%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.23"
$ws.Range("E10").Value = "  -1.07%  "
$ws.Range("E11").Value = "  -1.29%  "
$ws.Range("D12").Value = "4.094.13"
$ws.Range("E12").Value = "  +0.16%  "
$ws.Range("E13").Value = "  -0.42%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "29.06"
$ws.Range("E14").Value = "  +3.35%  "
$ws.Range("D15").Value = "66.872.37"
$ws.Range("E15").Value = "  +0.49%  "
$ws.Range("E16").Value = "  +0.35%  "
$ws.Range("D17").Value = "3.492.35"
$ws.Range("E17").Value = "  -0.60%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.28"
$ws.Range("E18").Value = "  -0.41%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.01"
$ws.Range("E19").Value = "  -0.23%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "394.63"
$ws.Range("E20").Value = "  +0.93%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.96"
$ws.Range("E21").Value = "  +0.80%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.13"
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  +0.23%  "
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000121"
$ws.Range("E25").Value = "  +0.24%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.18"
$ws.Range("E26").Value = "  +0.22%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.181"
$ws.Range("E27").Value = "  -0.43%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  +0.40%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.16"
$ws.Range("E29").Value = "  -2.24%  "
$ws.Range("E30").Value = "  -2.32%  "
$ws.Range("E31").Value = "  -0.09%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "23.70"
$ws.Range("E32").Value = "  +0.78%  "
$ws.Range("E33").Value = "  -0.87%  "
$ws.Range("E34").Value = "  +0.43%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "162.90"
$ws.Range("E35").Value = "  +0.89%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.876"
$ws.Range("E36").Value = "  -1.38%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.90"
$ws.Range("E37").Value = "  -0.81%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.97"
$ws.Range("E38").Value = "  +3.56%  "
$ws.Range("E39").Value = "  -0.04%  "
$ws.Range("E40").Value = "  -0.16%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "27.15"
$ws.Range("E41").Value = "  +1.24%  "
$ws.Range("D42").Value = "2.829.50"
$ws.Range("E42").Value = "  +2.20%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.23"
$ws.Range("E43").Value = "  -0.85%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "42.80"
$ws.Range("E44").Value = "  -0.60%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.55"
$ws.Range("E45").Value = "  +3.32%  "
$ws.Range("E46").Value = "  -2.79%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "335.93"
$ws.Range("E47").Value = "  -3.00%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "34.65"
$ws.Range("E48").Value = "  +2.69%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.08"
$ws.Range("E49").Value = "  -1.07%  "
$ws.Range("E50").Value = "  -1.77%  "
$ws.Range("B51").Value = "SuiNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.839"
$ws.Range("E51").Value = "  -5.28%  "
